$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: date moves from 44284 to 44277, values swap to the "150/11000/550" set
$ws.Range("D2").Value = 44277
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 11000
$ws.Range("L2").Value = 11000
$ws.Range("M2").Value = 11000
$ws.Range("P2").Value = 550

# Row 4: date moves from 44280 to 44284 (other values unchanged)
$ws.Range("D4").Value = 44284

# Row 5: date moves from 44277 to 44280, values swap to the "100/10000/500" set
$ws.Range("D5").Value = 44280
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 10000
$ws.Range("P5").Value = 500
